$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 144, shifting existing rows 144:252 down to 145:253.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new weekly record.
$ws.Range("A144").Value = 8
$ws.Range("B144").Value = "Terminal La Palmera de La Serena"
$ws.Range("C144").Value = "Coquimbo"
$ws.Range("D144").Value = 44651
$ws.Range("E144").Value = 4
$ws.Range("F144").Value = 100112012
$ws.Range("G144").Value = "Espinaca"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 2000
$ws.Range("K144").Value = 550
$ws.Range("L144").Value = 600
$ws.Range("M144").Value = 575
$ws.Range("N144").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O144").Value = "Provincia del Elquí"
$ws.Range("P144").Value = 1150
$ws.Range("Q144").Value = 0.5
$ws.Range("R144").Value = "Hortaliza"
